$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 18 (Leve Item ID 5471)
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -716

# Row 43 (Leve Item ID 5472)
$ws.Range("H43").Value = 3206.5
$ws.Range("I43").Value = 3542.111
$ws.Range("K43").Value = 3542.111
$ws.Range("M43").Value = -3473.111

# Row 64 (Leve Item ID 5506)
$ws.Range("H64").Value = 4518.8
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").Value = ""

# Row 67 (Leve Item ID 5506)
$ws.Range("H67").Value = 4518.8
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").Value = ""

# Row 74 (Leve Item ID 5507)
$ws.Range("H74").Value = 7148.6
$ws.Range("I74").Value = 7148.6
$ws.Range("K74").Value = 7148.6
$ws.Range("M74").Value = -6212.6

# Row 77 (Leve Item ID 5507)
$ws.Range("H77").Value = 7148.6
$ws.Range("I77").Value = 7148.6
$ws.Range("K77").Value = 35743
$ws.Range("M77").Value = -31063

# Row 88 (Leve Item ID 12608)
$ws.Range("H88").Value = 12684.75
$ws.Range("J88").Value = 12684.75
$ws.Range("L88").Value = 12684.75
$ws.Range("N88").Value = -13496.75

# Row 91 (Leve Item ID 12608)
$ws.Range("H91").Value = 12684.75
$ws.Range("J91").Value = 12684.75
$ws.Range("L91").Value = 12684.75
$ws.Range("N91").Value = -15492.75

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 3166.9722
$ws.Range("I138").Value = 2877.5
$ws.Range("J138").Value = 3278.3076
$ws.Range("K138").Value = 8632.5
$ws.Range("L138").Value = 9834.9228
$ws.Range("M138").Value = -3492.5
$ws.Range("N138").Value = -20114.9228


$ws = $wb.Worksheets.Item("ARM")

# Row 63 (Leve Item ID 12528)
$ws.Range("H63").Value = 995.8

# Row 66 (Leve Item ID 12528)
$ws.Range("H66").Value = 995.8

# Row 88 (Leve Item ID 12530)
$ws.Range("H88").Value = 1896.5
$ws.Range("I88").Value = 1971.6666
$ws.Range("J88").Value = 1835
$ws.Range("K88").Value = 1971.6666
$ws.Range("L88").Value = 1835
$ws.Range("M88").Value = -1565.6666
$ws.Range("N88").Value = -2647

# Row 91 (Leve Item ID 12530)
$ws.Range("H91").Value = 1896.5
$ws.Range("I91").Value = 1971.6666
$ws.Range("J91").Value = 1835
$ws.Range("K91").Value = 1971.6666
$ws.Range("L91").Value = 1835
$ws.Range("M91").Value = -567.6666
$ws.Range("N91").Value = -4643

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 4168972.5
$ws.Range("I132").Value = 4763826
$ws.Range("J132").Value = 4998.3335
$ws.Range("K132").Value = 14291478
$ws.Range("L132").Value = 14995.0005
$ws.Range("M132").Value = -14288948
$ws.Range("N132").Value = -20055.0005


$ws = $wb.Worksheets.Item("CRP")

# Row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 429.66666
$ws.Range("I22").Value = 243.5
$ws.Range("J22").Value = 578.6
$ws.Range("K22").Value = 243.5
$ws.Range("L22").Value = 578.6
$ws.Range("M22").Value = 106.5
$ws.Range("N22").Value = -1278.6

# Row 25 (Leve Item ID 1895)
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = ""

# Row 32 (Leve Item ID 2246)
$ws.Range("H32").Value = 2207.25
$ws.Range("I32").Value = 2207.25
$ws.Range("K32").Value = 2207.25
$ws.Range("M32").Value = -1891.25

# Row 41 (Leve Item ID 1917)
$ws.Range("H41").Value = 5450

# Row 51 (Leve Item ID 2039)
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""

# Row 60 (Leve Item ID 1937)
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").Value = ""

# Row 61 (Leve Item ID 2039)
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = ""


$ws = $wb.Worksheets.Item("CUL")

# Row 120 (Leve Item ID 27877)
$ws.Range("H120").Value = 17348.75
$ws.Range("I120").Value = 17348.75
$ws.Range("K120").Value = 52046.25
$ws.Range("M120").Value = -47208.25


$ws = $wb.Worksheets.Item("GSM")

# Row 3 (Leve Item ID 4091)
$ws.Range("H3").Value = 700
$ws.Range("I3").Value = 700
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 700
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -584
$ws.Range("N3").Value = ""

# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 3999.4
$ws.Range("I80").Value = 4249.5
$ws.Range("K80").Value = 4249.5
$ws.Range("M80").Value = -3251.5

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 3999.4
$ws.Range("I83").Value = 4249.5
$ws.Range("K83").Value = 21247.5
$ws.Range("M83").Value = -16255.5


$ws = $wb.Worksheets.Item("LTW")

# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 4280.4
$ws.Range("J7").Value = 4999.5
$ws.Range("L7").Value = 4999.5
$ws.Range("N7").Value = -5223.5

# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 2119.4666
$ws.Range("I16").Value = 1432.9166
$ws.Range("K16").Value = 1432.9166
$ws.Range("M16").Value = -1262.9166

# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 2950
$ws.Range("I22").Value = 3600
$ws.Range("K22").Value = 3600
$ws.Range("M22").Value = -3305

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 2950
$ws.Range("I27").Value = 3600
$ws.Range("K27").Value = 3600
$ws.Range("M27").Value = -3493

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 1357.4166
$ws.Range("I46").Value = 1371.8182
$ws.Range("J46").Value = 1199
$ws.Range("K46").Value = 1371.8182
$ws.Range("L46").Value = 1199
$ws.Range("M46").Value = -1183.8182
$ws.Range("N46").Value = -1575

# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 481.22223
$ws.Range("J55").Value = 650
$ws.Range("L55").Value = 650
$ws.Range("N55").Value = -996

# Row 56 (Leve Item ID 3668)
$ws.Range("H56").Value = 11499.5
$ws.Range("J56").Value = 11499.5
$ws.Range("L56").Value = 11499.5
$ws.Range("N56").Value = -12881.5

# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 2747.25
$ws.Range("I68").Value = 2731.3333
$ws.Range("J68").Value = 2795
$ws.Range("K68").Value = 2731.3333
$ws.Range("L68").Value = 2795
$ws.Range("M68").Value = -1982.3333
$ws.Range("N68").Value = -4293

# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 2747.25
$ws.Range("I71").Value = 2731.3333
$ws.Range("J71").Value = 2795
$ws.Range("K71").Value = 13656.6665
$ws.Range("L71").Value = 13975
$ws.Range("M71").Value = -9912.666499999999
$ws.Range("N71").Value = -21463

# Row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 1166.5
$ws.Range("I82").Value = 1190.6364
$ws.Range("J82").Value = 1033.75
$ws.Range("K82").Value = 1190.6364
$ws.Range("L82").Value = 1033.75
$ws.Range("M82").Value = -829.6364000000001
$ws.Range("N82").Value = -1755.75

# Row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 1166.5
$ws.Range("I85").Value = 1190.6364
$ws.Range("J85").Value = 1033.75
$ws.Range("K85").Value = 1190.6364
$ws.Range("L85").Value = 1033.75
$ws.Range("M85").Value = 57.36359999999991
$ws.Range("N85").Value = -3529.75

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 7383.9287
$ws.Range("I122").Value = 7079.636
$ws.Range("K122").Value = 21238.908
$ws.Range("M122").Value = -18788.908

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 4280.4
$ws.Range("J126").Value = 4999.5
$ws.Range("L126").Value = 14998.5
$ws.Range("N126").Value = -19938.5


$ws = $wb.Worksheets.Item("WVR")

# Row 62 (Leve Item ID 12589)
$ws.Range("H62").Value = 4033
$ws.Range("I62").Value = 2749.5
$ws.Range("J62").Value = 6600
$ws.Range("K62").Value = 2749.5
$ws.Range("L62").Value = 6600
$ws.Range("M62").Value = -2125.5
$ws.Range("N62").Value = -7848

# Row 65 (Leve Item ID 12589)
$ws.Range("H65").Value = 4033
$ws.Range("I65").Value = 2749.5
$ws.Range("J65").Value = 6600
$ws.Range("K65").Value = 13747.5
$ws.Range("L65").Value = 33000
$ws.Range("M65").Value = -10627.5
$ws.Range("N65").Value = -39240

